$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H106").Value = 5799.8
$ws.Range("I106").Value = 4333
$ws.Range("K106").Value = 4333
$ws.Range("M106").Value = -3702

$ws.Range("H111").Value = 40632.25
$ws.Range("I111").Value = 43722.57
$ws.Range("K111").Value = 131167.71
$ws.Range("M111").Value = -128100.71

$ws.Range("H125").Value = 3905.4
$ws.Range("J125").Value = 5700.5
$ws.Range("L125").Value = 51304.5
$ws.Range("N125").Value = -56224.5

$ws.Range("H138").Value = 3575.65
$ws.Range("I138").Value = 1449.4445
$ws.Range("J138").Value = 3785.934
$ws.Range("K138").Value = 4348.333500000001
$ws.Range("L138").Value = 11357.802
$ws.Range("M138").Value = 791.6664999999994
$ws.Range("N138").Value = -21637.802

$ws.Range("H141").Value = 6865
$ws.Range("I141").Value = 2999.8
$ws.Range("K141").Value = 8999.400000000001
$ws.Range("M141").Value = -3819.400000000001

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1250.138
$ws.Range("I2").Value = 1210.2307
$ws.Range("K2").Value = 1210.2307
$ws.Range("M2").Value = -1097.2307

$ws.Range("H116").Value = 1250.138
$ws.Range("I116").Value = 1210.2307
$ws.Range("K116").Value = 1210.2307
$ws.Range("M116").Value = 1083.7693

$ws.Range("H122").Value = 4006.4546
$ws.Range("I122").Value = 4137.0625
$ws.Range("J122").Value = 3658.1667
$ws.Range("K122").Value = 12411.1875
$ws.Range("L122").Value = 10974.5001
$ws.Range("M122").Value = -9961.1875
$ws.Range("N122").Value = -15874.5001

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H45").Value = 35000
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -36616

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H5").Value = 2754.3333
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""

$ws.Range("H99").Value = 4342.4814
$ws.Range("I99").Value = 3585.9443
$ws.Range("J99").Value = 5855.5557
$ws.Range("K99").Value = 3585.9443
$ws.Range("L99").Value = 5855.5557
$ws.Range("M99").Value = -2087.9443
$ws.Range("N99").Value = -8851.555700000001

$ws.Range("H122").Value = 4539.2104
$ws.Range("I122").Value = 4291.6665
$ws.Range("K122").Value = 12874.9995
$ws.Range("M122").Value = -10424.9995

$ws.Range("H126").Value = 4342.4814
$ws.Range("I126").Value = 3585.9443
$ws.Range("J126").Value = 5855.5557
$ws.Range("K126").Value = 10757.8329
$ws.Range("L126").Value = 17566.6671
$ws.Range("M126").Value = -8287.832900000001
$ws.Range("N126").Value = -22506.6671

$ws.Range("H134").Value = 2381.8696
$ws.Range("I134").Value = 821.2778
$ws.Range("K134").Value = 2463.8334
$ws.Range("M134").Value = 71.16660000000002

$ws.Range("H135").Value = 108250
$ws.Range("I135").Value = 106500
$ws.Range("J135").Value = 110000
$ws.Range("K135").Value = 106500
$ws.Range("L135").Value = 110000
$ws.Range("M135").Value = -101430
$ws.Range("N135").Value = -120140

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 161.57143
$ws.Range("I2").Value = 82.57143000000001
$ws.Range("K2").Value = 495.42858
$ws.Range("M2").Value = -382.42858

$ws.Range("H132").Value = 3635.1936
$ws.Range("I132").Value = 1993.25
$ws.Range("J132").Value = 3878.4443
$ws.Range("K132").Value = 17939.25
$ws.Range("L132").Value = 34905.9987
$ws.Range("M132").Value = -15409.25
$ws.Range("N132").Value = -39965.9987

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H9").Value = 268.66666
$ws.Range("I9").Value = 268.66666
$ws.Range("K9").Value = 268.66666
$ws.Range("M9").Value = -98.66665999999998

$ws.Range("H11").Value = 10750000
$ws.Range("I11").Value = 12187500
$ws.Range("K11").Value = 12187500
$ws.Range("M11").Value = -12187361

$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""

$ws.Range("H14").Value = 11500000
$ws.Range("I14").Value = 11500000
$ws.Range("K14").Value = 11500000
$ws.Range("M14").Value = -11499832

$ws.Range("H126").Value = 3313.6155
$ws.Range("I126").Value = 2947.5
$ws.Range("J126").Value = 4534
$ws.Range("K126").Value = 8842.5
$ws.Range("L126").Value = 13602
$ws.Range("M126").Value = -6372.5
$ws.Range("N126").Value = -18542

$ws.Range("H133").Value = 69518.664
$ws.Range("J133").Value = 69518.664
$ws.Range("L133").Value = 69518.664
$ws.Range("N133").Value = -79638.664

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H3").Value = 5756.6
$ws.Range("J3").Value = 5756.6
$ws.Range("L3").Value = 5756.6
$ws.Range("N3").Value = -5980.6

$ws.Range("H15").Value = 5756.6
$ws.Range("J15").Value = 5756.6
$ws.Range("L15").Value = 5756.6
$ws.Range("N15").Value = -6096.6

$ws.Range("H40").Value = 6579.08
$ws.Range("I40").Value = 5824.8237
$ws.Range("K40").Value = 5824.8237
$ws.Range("M40").Value = -5688.8237

$ws.Range("H61").Value = 89332.164
$ws.Range("I61").Value = 173666.33
$ws.Range("K61").Value = 173666.33
$ws.Range("M61").Value = -173464.33

$ws.Range("H104").Value = 88333.336
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 88333.336
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 88333.336
$ws.Range("M104").Value = ""
$ws.Range("N104").Value = -95321.336

$ws.Range("H113").Value = 89332.164
$ws.Range("I113").Value = 173666.33
$ws.Range("K113").Value = 173666.33
$ws.Range("M113").Value = -171496.33

$ws.Range("H132").Value = 3349.1892
$ws.Range("I132").Value = 3258.4092
$ws.Range("K132").Value = 9775.2276
$ws.Range("M132").Value = -7245.2276

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""

$ws.Range("H38").Value = 8499.5
$ws.Range("I38").Value = 8499
$ws.Range("K38").Value = 8499
$ws.Range("M38").Value = -8026

$ws.Range("H39").Value = 29543.75
$ws.Range("J39").Value = 29543.75
$ws.Range("L39").Value = 29543.75
$ws.Range("N39").Value = -30369.75

$ws.Range("H122").Value = 2329.4138
$ws.Range("I122").Value = 1841.9524
$ws.Range("J122").Value = 3609
$ws.Range("K122").Value = 5525.857199999999
$ws.Range("L122").Value = 10827
$ws.Range("M122").Value = -3075.857199999999
$ws.Range("N122").Value = -15727

$ws.Range("H132").Value = 995.4857
$ws.Range("I132").Value = 1035.5
$ws.Range("J132").Value = 879.8889
$ws.Range("K132").Value = 3106.5
$ws.Range("L132").Value = 2639.6667
$ws.Range("M132").Value = -576.5
$ws.Range("N132").Value = -7699.6667

$ws.Range("H135").Value = 60514.91
$ws.Range("J135").Value = 60514.91
$ws.Range("L135").Value = 60514.91
$ws.Range("N135").Value = -70654.91
